$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF), per the commit
# "repull data, push all data, mean calculation"
$updates = @{
    2  = -1
    5  = 0
    6  = 3
    7  = 0
    8  = 3
    15 = 3
    16 = 2
    19 = 0
    22 = 4
    31 = -2
    33 = 0
    35 = 2
    37 = 1
    43 = -2
    46 = 2
    49 = -4
    52 = 1
    57 = 3
    58 = 4
    61 = 5
    67 = 4
    71 = -1
    72 = -3
    73 = 5
    74 = 3
    75 = 1
    77 = 8
    78 = 1
    79 = 4
    80 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
